$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace mock PII "Name" column values with cartoon mock names (H2:H7)
$ws.Range("H2").Value = "Mickey Mouse"
$ws.Range("H3").Value = "Mickey Mouse"
$ws.Range("H4").Value = "Donald Duck"
$ws.Range("H5").Value = "Donald Duck"
$ws.Range("H6").Value = "Jane Doe"
$ws.Range("H7").Value = "John Doe"

# Update sheet view: scroll so column B is leftmost, select N2:N7
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("N2:N7").Select()
